$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 757
$ws.Range("E2").Value = 108
$ws.Range("F2").Value = -41
$ws.Range("G2").Value = 82
$ws.Range("H2").Value = 51
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = -15
$ws.Range("K2").Value = 1035
$ws.Range("L2").Value = 308
$ws.Range("M2").Value = 727
$ws.Range("N2").Value = 709
$ws.Range("O2").Value = 18
$ws.Range("P2").Value = 494
$ws.Range("Q2").Value = 418
$ws.Range("R2").Value = -346
$ws.Range("S2").Value = 12
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 418
$ws.Range("V2").Value = 55
$ws.Range("W2").Value = 14.32
$ws.Range("X2").Value = 6.73
$ws.Range("Y2").Value = 12.4
$ws.Range("Z2").Value = 5.51
$ws.Range("AA2").Value = 42.28
$ws.Range("AB2").Value = 109.89
$ws.Range("AC2").Value = 160
$ws.Range("AD2").Value = 19.68
$ws.Range("AE2").Value = 1650
$ws.Range("AF2").Value = 1.91
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 43668187

# Row 3
$ws.Range("D3").Value = 1062
$ws.Range("E3").Value = 74
$ws.Range("F3").Value = 74
$ws.Range("G3").Value = 65
$ws.Range("H3").Value = 38
$ws.Range("I3").Value = 51
$ws.Range("J3").Value = -13
$ws.Range("K3").Value = 2557
$ws.Range("L3").Value = 406
$ws.Range("M3").Value = 2150
$ws.Range("N3").Value = 2134
$ws.Range("O3").Value = 17
$ws.Range("P3").Value = 717
$ws.Range("Q3").Value = 346
$ws.Range("R3").Value = -208
$ws.Range("S3").Value = -96
$ws.Range("T3").Value = 8
$ws.Range("U3").Value = 339
$ws.Range("V3").Value = 75
$ws.Range("W3").Value = 6.94
$ws.Range("X3").Value = 3.6
$ws.Range("Y3").Value = 3.61
$ws.Range("Z3").Value = 2.13
$ws.Range("AA3").Value = 18.9
$ws.Range("AB3").Value = 206.05
$ws.Range("AC3").Value = 42
$ws.Range("AD3").Value = 49.22
$ws.Range("AE3").Value = 1518
$ws.Range("AF3").Value = 1.35
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 143309751

# Row 4
$ws.Range("D4").Value = 1094
$ws.Range("E4").Value = 112
$ws.Range("F4").Value = 112
$ws.Range("G4").Value = 62
$ws.Range("H4").Value = 42
$ws.Range("I4").Value = 44
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 2566
$ws.Range("L4").Value = 351
$ws.Range("M4").Value = 2216
$ws.Range("N4").Value = 2200
$ws.Range("O4").Value = 15
$ws.Range("P4").Value = 718
$ws.Range("Q4").Value = 422
$ws.Range("R4").Value = -383
$ws.Range("S4").Value = -28
$ws.Range("T4").Value = 7
$ws.Range("U4").Value = 415
$ws.Range("V4").Value = 40
$ws.Range("W4").Value = 10.26
$ws.Range("X4").Value = 3.86
$ws.Range("Y4").Value = 2.01
$ws.Range("Z4").Value = 1.65
$ws.Range("AA4").Value = 15.83
$ws.Range("AB4").Value = 212.31
$ws.Range("AC4").Value = 30
$ws.Range("AD4").Value = 61.37
$ws.Range("AE4").Value = 1561
$ws.Range("AF4").Value = 1.19
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 143556056

# Row 5
$ws.Range("D5").Value = 1268
$ws.Range("E5").Value = 153
$ws.Range("F5").Value = 153
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = 26
$ws.Range("I5").Value = 25
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2551
$ws.Range("L5").Value = 274
$ws.Range("M5").Value = 2277
$ws.Range("N5").Value = 2261
$ws.Range("O5").Value = 16
$ws.Range("P5").Value = 726
$ws.Range("Q5").Value = 618
$ws.Range("R5").Value = -509
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 28
$ws.Range("U5").Value = 590
$ws.Range("V5").Value = 20
$ws.Range("W5").Value = 12.1
$ws.Range("X5").Value = 2.08
$ws.Range("Y5").Value = 1.12
$ws.Range("Z5").Value = 1.03
$ws.Range("AA5").Value = 12.01
$ws.Range("AB5").Value = 216.53
$ws.Range("AC5").Value = 17
$ws.Range("AD5").Value = 146.39
$ws.Range("AE5").Value = 1586
$ws.Range("AF5").Value = 1.6
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 145153497

# Row 6
$ws.Range("D6").Value = 1094
$ws.Range("E6").Value = 99
$ws.Range("F6").Value = 99
$ws.Range("G6").Value = -64
$ws.Range("H6").Value = -99
$ws.Range("I6").Value = -101
$ws.Range("K6").Value = 2313
$ws.Range("L6").Value = 267
$ws.Range("M6").Value = 2046
$ws.Range("N6").Value = 2028
$ws.Range("P6").Value = 731
$ws.Range("Q6").Value = 405
$ws.Range("R6").Value = -571
$ws.Range("S6").Value = -149
$ws.Range("T6").Value = 11
$ws.Range("U6").Value = 394
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 9.07
$ws.Range("X6").Value = -9.039999999999999
$ws.Range("Y6").Value = -4.71
$ws.Range("Z6").Value = -4.07
$ws.Range("AA6").Value = 13.05
$ws.Range("AB6").Value = 182.2
$ws.Range("AC6").Value = -69
$ws.Range("AD6").Value = -26.89
$ws.Range("AE6").Value = 1412
$ws.Range("AF6").Value = 1.32
$ws.Range("AG6").Value = 104
$ws.Range("AH6").Value = 5.59
$ws.Range("AI6").Value = -147.94
$ws.Range("AJ6").Value = 146235748

# Row 7
$ws.Range("D7").Value = 937
$ws.Range("E7").Value = -42
$ws.Range("G7").Value = -25
$ws.Range("H7").Value = -30
$ws.Range("I7").Value = -27
$ws.Range("K7").Value = 2270
$ws.Range("L7").Value = 250
$ws.Range("M7").Value = 2020
$ws.Range("N7").Value = 2000
$ws.Range("P7").Value = 730
$ws.Range("Q7").Value = 500
$ws.Range("R7").Value = -30
$ws.Range("S7").Value = -10
$ws.Range("T7").Value = 0
$ws.Range("W7").Value = -4.48
$ws.Range("X7").Value = -3.2
$ws.Range("Y7").Value = -1.34
$ws.Range("Z7").Value = -1.31
$ws.Range("AA7").Value = 12.38
$ws.Range("AC7").Value = -18
$ws.Range("AD7").Value = -87.47
$ws.Range("AE7").Value = 1392
$ws.Range("AF7").Value = 1.16
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("U7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 1022
$ws.Range("E8").Value = 100
$ws.Range("G8").Value = 131
$ws.Range("H8").Value = 120
$ws.Range("I8").Value = 123
$ws.Range("K8").Value = 2410
$ws.Range("L8").Value = 270
$ws.Range("M8").Value = 2140
$ws.Range("N8").Value = 2120
$ws.Range("P8").Value = 730
$ws.Range("Q8").Value = 480
$ws.Range("R8").Value = -40
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("W8").Value = 9.789999999999999
$ws.Range("X8").Value = 11.74
$ws.Range("Y8").Value = 5.97
$ws.Range("Z8").Value = 5.13
$ws.Range("AA8").Value = 12.62
$ws.Range("AC8").Value = 84
$ws.Range("AD8").Value = 19.2
$ws.Range("AE8").Value = 1476
$ws.Range("AF8").Value = 1.09
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("U8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 1100
$ws.Range("E9").Value = 120
$ws.Range("G9").Value = 160
$ws.Range("H9").Value = 150
$ws.Range("I9").Value = 150
$ws.Range("K9").Value = 2580
$ws.Range("L9").Value = 290
$ws.Range("M9").Value = 2290
$ws.Range("N9").Value = 2270
$ws.Range("P9").Value = 730
$ws.Range("Q9").Value = 340
$ws.Range("R9").Value = -40
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("W9").Value = 10.91
$ws.Range("X9").Value = 13.64
$ws.Range("Y9").Value = 6.83
$ws.Range("Z9").Value = 6.01
$ws.Range("AA9").Value = 12.66
$ws.Range("AC9").Value = 103
$ws.Range("AD9").Value = 15.74
$ws.Range("AE9").Value = 1580
$ws.Range("AF9").Value = 1.02
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("U9").ClearContents()
$ws.Range("AI9").ClearContents()
